$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Webster's Method row (row 10): update simulation vehicle count and travel time
$ws.Range("C10").Value = 716
$ws.Range("D10").Value = "127.6 saniye"

# Update the active selection to reflect where the edit was made
$ws.Range("D10").Select()
